$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 290. This shifts existing rows 290-389
# down to 292-391, carrying their values/formats with them.
$ws.Rows.Item(290).Resize(2).Insert()

# New row 290: Primera, 2022-xx-xx (serial 44988)
$ws.Range("A290").Value = 7
$ws.Range("B290").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C290").Value = "Ñuble"
$ws.Range("D290").Value = 44988
$ws.Range("D290").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E290").Value = 16
$ws.Range("F290").Value = 100112009
$ws.Range("G290").Value = "Acelga"
$ws.Range("H290").Value = "Sin especificar"
$ws.Range("I290").Value = "Primera"
$ws.Range("J290").Value = 400
$ws.Range("K290").Value = 600
$ws.Range("L290").Value = 700
$ws.Range("M290").Value = 650
$ws.Range("N290").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O290").Value = "Provincia de Diguillín"
$ws.Range("P290").Value = 650
$ws.Range("Q290").Value = 1
$ws.Range("R290").Value = "Hortaliza"

# New row 291: Segunda, same date
$ws.Range("A291").Value = 7
$ws.Range("B291").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C291").Value = "Ñuble"
$ws.Range("D291").Value = 44988
$ws.Range("D291").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E291").Value = 16
$ws.Range("F291").Value = 100112009
$ws.Range("G291").Value = "Acelga"
$ws.Range("H291").Value = "Sin especificar"
$ws.Range("I291").Value = "Segunda"
$ws.Range("J291").Value = 300
$ws.Range("K291").Value = 500
$ws.Range("L291").Value = 500
$ws.Range("M291").Value = 500
$ws.Range("N291").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O291").Value = "Provincia de Diguillín"
$ws.Range("P291").Value = 500
$ws.Range("Q291").Value = 1
$ws.Range("R291").Value = "Hortaliza"
